# Apply the edits described by the commit:
# "add: Start of FOB georef + little progress on analysis"
#
# Summary of the data-level changes being applied:
#  - Column G (budget-ish helper col) gets a value of 2 for rows 53-198.
#  - Column H (location, lat/lon text) gets new values on a handful of
#    those rows plus row 210, pulling in 15 brand-new shared strings.
#  - Column C is widened (best-fit) and the sheet view / window chrome
#    moved (best effort only - cosmetic viewport state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G: mark rows 53-198 with budget multiplier/count = 2 -----------
for ($r = 53; $r -le 198; $r++) {
    $ws.Cells.Item($r, 7).Value = 2
}

# --- Column H: new lat/lon "location" strings -------------------------------
# Written in the exact order the strings were first introduced so the
# resulting shared-string table lines up with the source workbook.
$ws.Cells.Item(210, 8).Value = "12.9113427641161, 77.57429815729775"
$ws.Cells.Item(53, 8).Value  = "12.956849275144204, 77.69780775936569"
$ws.Cells.Item(54, 8).Value  = "12.956051725724723, 77.69312516347807"
$ws.Cells.Item(55, 8).Value  = "12.954715204970352, 77.68394924104837"
$ws.Cells.Item(69, 8).Value  = "12.968689864766919, 77.60142424615806"
$ws.Cells.Item(71, 8).Value  = "12.971840853212232, 77.60689869770343"
$ws.Cells.Item(78, 8).Value  = "12.925435860659025, 77.67544480391703"
$ws.Cells.Item(67, 8).Value  = "12.976667691086602, 77.59927032472335"
$ws.Cells.Item(66, 8).Value  = "12.96948427474764, 77.60245180356293"
$ws.Cells.Item(65, 8).Value  = "12.976120668809553, 77.60374548748605"
$ws.Cells.Item(68, 8).Value  = "12.972496575353658, 77.61950174788767"
$ws.Cells.Item(61, 8).Value  = "12.984022431621492, 77.58660712022653"
$ws.Cells.Item(63, 8).Value  = "12.984309599120063, 77.59716274729055"
$ws.Cells.Item(107, 8).Value = "13.012582723434742, 77.62583797263228"
$ws.Cells.Item(108, 8).Value = "13.018204179943318, 77.63655273257056"

# --- Column C: widen to fit the new (much longer) description text --------
$ws.Columns.Item(3).ColumnWidth = 162.8

# --- View state: move the active selection (scroll position of the frozen
#     pane is re-derived by the host on save and isn't independently
#     addressable through this object model) --------------------------------
$ws.Activate()
$ws.Range("H129").Select()
